$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.0725292808699395, 1.452525492729129, 3.943704826413877, 1.985876337140326, 2.004298684434776, 51)
    3  = @(0.970128790070132, 1.91197192143461, 11.00797031802001, 3.317826143428858, 3.205037560462687, 50)
    4  = @(1.027439519176732, 2.670308812046868, 20.8551199302896, 4.566740624372004, 4.495773270230817, 49)
    5  = @(1.040770019679559, 2.551126684817124, 21.80091093632453, 4.669144561514939, 4.599838150987386, 48)
    6  = @(0.6404470419866763, 2.649364789166133, 22.40911179286584, 4.733826337421541, 4.741010144670577, 47)
    7  = @(0.6619285698089534, 2.518113852960962, 23.16329070585441, 4.812825646733363, 4.83108003263652, 38)
    8  = @(0.6925657535171938, 2.61980777971489, 24.31016614921836, 4.930534063285474, 4.948987541160202, 37)
    9  = @(0.7380237456049207, 3.434955442878508, 40.91359093038702, 6.396373263841551, 6.518710873632213, 20)
    10 = @(-1.030652359722206, 3.294535520879184, 28.3461002497294, 5.324105582135783, 5.436681950116494, 13)
    11 = @(-0.1548825442147432, 2.972445703758152, 9.566860607955785, 3.093034207369163, 3.45377909062185, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
